$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header row: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Convert the data range into an Excel Table ("Table1") while preserving the
#    original header-row formatting (stash it on an unused scratch row, clear
#    the header formatting so the table creation does not "capture" a header
#    dxf, add the table, then restore the stashed formatting and discard the
#    scratch row again).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A1000:U1000")

$headerRange.Copy() | Out-Null
$scratchRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$headerRange.ClearFormats() | Out-Null

$tableRange = $ws.Range("A1:U78")
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

$scratchRange.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null

$scratchRange.EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split after row 1) and leave the selection in the
#    (bottom-left) scrollable pane, right below the frozen header.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
